$d = $word.ActiveDocument

$replacements = @(
    @{old="85×89="; new="78×23="},
    @{old="36×71="; new="77×74="},
    @{old="52×37="; new="33×96="},
    @{old="30×19="; new="53×73="},
    @{old="15×44="; new="62×60="},
    @{old="48×34="; new="83×82="},
    @{old="16×51="; new="55×50="},
    @{old="52×68="; new="87×26="},
    @{old="61×70="; new="79×62="},
    @{old="77×21="; new="35×20="},
    @{old="46×33="; new="12×21="},
    @{old="31×73="; new="28×12="},
    @{old="43×99="; new="19×91="},
    @{old="81×94="; new="25×75="},
    @{old="16×61="; new="72×88="},
    @{old="29×78="; new="28×91="},
    @{old="28×21="; new="33×41="},
    @{old="45×51="; new="89×63="},
    @{old="99×88="; new="22×41="},
    @{old="86×50="; new="26×62="},
    @{old="63×33="; new="37×80="},
    @{old="81×23="; new="20×70="},
    @{old="33×58="; new="74×44="},
    @{old="31×83="; new="75×12="},
    @{old="11×96="; new="82×65="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
